$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '27.820.58'
$ws.Range("E2").Value = '  +0.38%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.908.83'
$ws.Range("E3").Value = '  +0.78%  '

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9996'
$ws.Range("E4").Value = '  -0.08%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '312.91'
$ws.Range("E5").Value = '  +0.38%  '

# Row 6
$ws.Range("E6").Value = '  +0.02%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5241'
$ws.Range("E7").Value = '  +7.33%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3788'
$ws.Range("E8").Value = '  -0.17%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07250'
$ws.Range("E9").Value = '  -1.04%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '21.36'
$ws.Range("E10").Value = '  +3.83%  '

# Row 11
$ws.Range("E11").Value = '  -0.64%  '

# Row 12
$ws.Range("B12").Value = 'TRON'
$ws.Range("C12").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07662'
$ws.Range("E12").Value = '  +0.16%  '

# Row 13
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.908.17'
$ws.Range("E13").Value = '  +0.71%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.457'
$ws.Range("E14").Value = '  -0.34%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '92.10'
$ws.Range("E15").Value = '  +0.85%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.9996'
$ws.Range("E16").Value = '  -0.20%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008729'
$ws.Range("E17").Value = '  -0.41%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.000'
$ws.Range("E18").Value = '  -0.02%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '27.849.73'
$ws.Range("E19").Value = '  +0.20%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.55'
$ws.Range("E20").Value = '  +0.62%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.162'
$ws.Range("E21").Value = '  +0.92%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.182.69'
$ws.Range("E22").Value = '  +2.03%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.86'
$ws.Range("E23").Value = '  +1.05%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.654'
$ws.Range("E24").Value = '  +0.43%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '154.01'
$ws.Range("E25").Value = '  +0.00%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.870'
$ws.Range("E26").Value = '  -0.69%  '

# Row 27
$ws.Range("B27").Value = 'LidoDAOToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.172'
$ws.Range("E27").Value = '  +0.80%  '

# Row 28
$ws.Range("B28").Value = 'EthereumClassic'
$ws.Range("C28").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '18.37'
$ws.Range("E28").Value = '  +0.06%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '115.01'
$ws.Range("E29").Value = '  -0.22%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.864'
$ws.Range("E30").Value = '  -0.09%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.09096'
$ws.Range("E31").Value = '  +2.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.864'
$ws.Range("E32").Value = '  +4.97%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.184'
$ws.Range("E33").Value = '  -0.47%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.235'
$ws.Range("E34").Value = '  +0.64%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7808'
$ws.Range("E35").Value = '  +1.62%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02098'
$ws.Range("E36").Value = '  +2.96%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.613'
$ws.Range("E37").Value = '  +1.81%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.072'
$ws.Range("E38").Value = '  +3.10%  '

# Row 39
$ws.Range("E39").Value = '  +2.46%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.094'
$ws.Range("E40").Value = '  -0.09%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.05284'
$ws.Range("E41").Value = '  +0.12%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.730'
$ws.Range("E42").Value = '  -2.36%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '115.88'
$ws.Range("E43").Value = '  +2.83%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.577'
$ws.Range("E44").Value = '  +0.75%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1523'
$ws.Range("E45").Value = '  +0.18%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4821'
$ws.Range("E46").Value = '  +0.66%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.50'
$ws.Range("E47").Value = '  -0.75%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("E48").Value = '  +0.08%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.623'
$ws.Range("E49").Value = '  -0.94%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '67.03'
$ws.Range("E50").Value = '  -0.67%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.05989'
$ws.Range("E51").Value = '  -0.99%  '

